# Fruta / hortaliza, semanal
#
# A new weekly price observation is inserted as row 109 on the single
# worksheet (Hortaliza, Feria Lagunitas de Puerto Montt - Ají). All the
# existing rows from 109..144 shift down by one (to 110..145) and the
# sheet's used range grows from A1:R144 to A1:R145.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 109..144 down to 110..145, leaving a blank row 109 behind.
$ws.Rows("109").Insert()

# Populate the newly inserted row 109 with the new observation.
$ws.Range("A109").Value = 4
$ws.Range("B109").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C109").Value = 'Los Lagos'
$ws.Range("D109").Value = 44463
$ws.Range("E109").Value = 10
$ws.Range("F109").Value = 100112021
$ws.Range("G109").Value = 'Ají'
$ws.Range("H109").Value = 'Inferno'
$ws.Range("I109").Value = 'Primera'
$ws.Range("J109").Value = 140
$ws.Range("K109").Value = 48000
$ws.Range("L109").Value = 50000
$ws.Range("M109").Value = 49000
$ws.Range("N109").Value = '$/caja 12 kilos'
$ws.Range("O109").Value = 'Región de Arica y Parinacota'
$ws.Range("P109").Value = 4083
$ws.Range("Q109").Value = 12
$ws.Range("R109").Value = 'Hortaliza'
